$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 2.25
$ws.Range("H2").Value = 3.4
$ws.Range("I2").Value = 3.1
$ws.Range("J2").Value = 3
$ws.Range("L2").Value = 3.75
$ws.Range("N2").Value = 10
$ws.Range("O2").Value = 1.3
$ws.Range("P2").Value = 3.4
$ws.Range("Q2").Value = 2.03
$ws.Range("R2").Value = 1.78
$ws.Range("U2").Value = 1.8
$ws.Range("V2").Value = 1.91
$ws.Range("AH2").Value = 9.5
$ws.Range("I3").Value = 3.3
$ws.Range("M3").Value = 1.07
$ws.Range("N3").Value = 9
$ws.Range("X3").Value = 9.5
$ws.Range("AX3").Value = 21
$ws.Range("G4").Value = 2.35
$ws.Range("J4").Value = 3.1
$ws.Range("K4").Value = 2.05
$ws.Range("L4").Value = 3.6
$ws.Range("U4").Value = 1.83
$ws.Range("V4").Value = 1.83
$ws.Range("Z4").Value = 23
$ws.Range("AA4").Value = 21
$ws.Range("AC4").Value = 9
$ws.Range("AG4").Value = 301
$ws.Range("AH4").Value = 8.5
$ws.Range("AK4").Value = 29
$ws.Range("AS4").Value = 201
$ws.Range("G6").Value = 1.42
$ws.Range("H6").Value = 4.2
$ws.Range("I6").Value = 8
$ws.Range("J6").Value = 1.95
$ws.Range("K6").Value = 2.3
$ws.Range("L6").Value = 7.5
$ws.Range("M6").Value = 1.05
$ws.Range("N6").Value = 11
$ws.Range("O6").Value = 1.3
$ws.Range("P6").Value = 3.5
$ws.Range("Q6").Value = 1.98
$ws.Range("R6").Value = 1.88
$ws.Range("U6").Value = 2.2
$ws.Range("V6").Value = 1.62
$ws.Range("Z6").Value = 9
$ws.Range("AB6").Value = 29
$ws.Range("AC6").Value = 9.5
$ws.Range("AK6").Value = 101
$ws.Range("AM6").Value = 67
$ws.Range("AQ6").Value = 21
$ws.Range("AU6").Value = 10
$ws.Range("AW6").Value = 8.5
$ws.Range("BC6").Value = 126
$ws.Range("Q7").Value = 2
$ws.Range("R7").Value = 1.8
$ws.Range("G8").Value = 2.2
$ws.Range("I8").Value = 3.2
$ws.Range("Q8").Value = 2
$ws.Range("R8").Value = 1.85
$ws.Range("BD9").Value = 126
$ws.Range("Q10").Value = 1.48
$ws.Range("R10").Value = 2.6
$ws.Range("G11").Value = 2.8
$ws.Range("I11").Value = 2.3
$ws.Range("J11").Value = 3.4
$ws.Range("N11").Value = 12
$ws.Range("AD11").Value = 6.5
$ws.Range("AK11").Value = 23
$ws.Range("AL11").Value = 19
$ws.Range("M12").Value = 1.05
$ws.Range("N12").Value = 11
$ws.Range("Q12").Value = 1.93
$ws.Range("R12").Value = 1.93
$ws.Range("M13").Value = 1.05
$ws.Range("N13").Value = 11
$ws.Range("O13").Value = 1.29
$ws.Range("P13").Value = 3.5
$ws.Range("Q13").Value = 1.9
$ws.Range("R13").Value = 1.95
$ws.Range("O17").Value = 1.33
$ws.Range("P17").Value = 3.25
$ws.Range("G18").Value = 1.62
$ws.Range("H18").Value = 4.1
$ws.Range("I18").Value = 4.75
$ws.Range("J18").Value = 2.2
$ws.Range("W18").Value = 7.5
$ws.Range("X18").Value = 8
$ws.Range("Z18").Value = 12
$ws.Range("AD18").Value = 8
$ws.Range("AL18").Value = 41
$ws.Range("AN18").Value = 3.6
$ws.Range("AO18").Value = 8
$ws.Range("AW18").Value = 7
$ws.Range("G20").Value = 5
$ws.Range("H20").Value = 4.1
$ws.Range("I20").Value = 1.6
$ws.Range("Q20").Value = 1.65
$ws.Range("R20").Value = 2.2
$ws.Range("Y20").Value = 15
$ws.Range("AA20").Value = 34
$ws.Range("AD20").Value = 8
$ws.Range("AH20").Value = 8.5
$ws.Range("N21").Value = 15
$ws.Range("O21").Value = 1.18
$ws.Range("P21").Value = 4.5
$ws.Range("Q21").Value = 1.6
$ws.Range("R21").Value = 2.3
$ws.Range("AL21").Value = 51
$ws.Range("AN21").Value = 3.5
$ws.Range("G22").Value = 2.55
$ws.Range("I22").Value = 2.8
$ws.Range("L22").Value = 3.4
$ws.Range("AI22").Value = 15
$ws.Range("AJ22").Value = 11
$ws.Range("AK22").Value = 29
$ws.Range("AN22").Value = 4.5
$ws.Range("AO22").Value = 13
$ws.Range("AQ22").Value = 41
